$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Fill in the new / moved entries.
# The order in which brand-new text values are assigned matters: Excel
# appends each new distinct string to the shared string table the first
# time it is used, so we replay the values in the same order the workbook
# author originally typed them (D17, D21, R19, U19, D18) to keep the
# shared-string indices identical to the target file.
# ---------------------------------------------------------------------------

# Row 17 : 2022-01-27, "12.00  - 19.30", 7.5 hours (new)
$ws.Range("C17").Value = 44588
$ws.Range("D17").Value = "12.00  - 19.30"
$ws.Range("F17").Value = 7.5

# Row 21 : 2022-02-02, "12.00 - 21.00", no hours yet (new, in progress)
$ws.Range("C21").Value = 44594
$ws.Range("D21").Value = "12.00 - 21.00"

# Row 19 (Learning goal 2 side) : 2022-01-28, "Menu - level transition" (new)
$ws.Range("Q19").Value = 44589
$ws.Range("R19").Value = "11.00 - 13.30"
$ws.Range("T19").Value = 2.5
$ws.Range("U19").Value = "Menu - level transition"

# Row 18 : 2022-01-28, "14.00 - 17.00", 3 hours (new)
$ws.Range("C18").Value = 44589
$ws.Range("D18").Value = "14.00 - 17.00"
$ws.Range("F18").Value = 3

# ---------------------------------------------------------------------------
# Learning goal 1 columns (C=Date, D=From-to, F=Hours) - remaining rows
# ---------------------------------------------------------------------------

# Row 16 : 2022-01-26, "12.00 - 20.00", 8 hours (previously on row 18)
$ws.Range("C16").Value = 44587
$ws.Range("D16").Value = "12.00 - 20.00"
$ws.Range("F16").Value = 8

# Row 19 : 2022-01-31, "13.00 - 17.00", 4 hours (previously on row 21, date shifted)
$ws.Range("C19").Value = 44592
$ws.Range("D19").Value = "13.00 - 17.00"
$ws.Range("F19").Value = 4

# Row 20 : 2022-02-01, "13.00 - 17.00", 4 hours (previously on row 22, date shifted)
$ws.Range("C20").Value = 44593
$ws.Range("D20").Value = "13.00 - 17.00"
$ws.Range("F20").Value = 4

# Row 21 no longer has an hours value (task still in progress)
$ws.Range("F21").ClearContents()

# Row 22 : no data yet, only keep the date number format
$ws.Range("C22:F22").ClearContents()

# Row 23 : new blank row, keep date format only
$ws.Range("C23").NumberFormat = "d-mmm"

# ---------------------------------------------------------------------------
# Learning goal 2 columns (Q=Date, R=From-to, T=Hours, U=Activity) - rest
# ---------------------------------------------------------------------------

# Row 15 : 2022-01-09 entry (shifted up from row 16)
$ws.Range("Q15").Value = 44570
$ws.Range("R15").Value = "21.45 - 22.45"
$ws.Range("T15").Value = 1
$ws.Range("U15").Value = "Implemented game over trigger"

# Row 16 : 2022-01-10 entry (shifted up from row 17)
$ws.Range("Q16").Value = 44571
$ws.Range("R16").Value = "10.30 - 13.00"
$ws.Range("T16").Value = 2.5
$ws.Range("U16").Value = "Added lava"

# Row 17 : 2022-01-12 entry (shifted up from row 18)
$ws.Range("Q17").Value = 44573
$ws.Range("R17").Value = "9.00 - 12.30"
$ws.Range("T17").Value = 3.5
$ws.Range("U17").Value = "Worked on level"

# Row 18 : 2022-01-13 entry (shifted up from row 19, values unchanged)
$ws.Range("Q18").Value = 44574
$ws.Range("R18").Value = "8.30 - 12.00"
$ws.Range("T18").Value = 3.5
$ws.Range("U18").Value = "Worked on level"

# Row 20 : 2022-01-31 entry (shifted up from row 24, no activity text)
$ws.Range("Q20").Value = 44592
$ws.Range("R20").Value = "9.00 - 12.00"
$ws.Range("T20").Value = 3
$ws.Range("U20").ClearContents()

# Row 21 : 2022-02-01 entry (shifted up from row 25, no activity text)
$ws.Range("Q21").Value = 44593
$ws.Range("R21").Value = "9.00 - 12.00"
$ws.Range("T21").Value = 3
$ws.Range("U21").ClearContents()

# Row 22 : 2022-02-02 entry (new, no activity text)
$ws.Range("Q22").Value = 44594
$ws.Range("R22").Value = "9.00 - 12.00"
$ws.Range("T22").Value = 3
$ws.Range("U22").ClearContents()

# Row 23 : no longer has data, just keep the date format
$ws.Range("Q23:T23").ClearContents()

# Row 25 / 26 : blank rows with just date number format
$ws.Range("Q25:T25").ClearContents()
$ws.Range("Q25").NumberFormat = "d-mmm"
$ws.Range("Q26").NumberFormat = "d-mmm"

# ---------------------------------------------------------------------------
# Totals
# ---------------------------------------------------------------------------
$ws.Range("F24").Formula = "=SUM(F3:F23)"
$ws.Range("Q24:T24").ClearContents()
$ws.Range("Q24").NumberFormat = "d-mmm"
$ws.Range("T24").Formula = "=SUM(T3:T22)"

# Remove the old row 27 (its total now lives in T24)
$ws.Rows(27).ClearContents()

$excel.Calculate()

# ---------------------------------------------------------------------------
# View state
# ---------------------------------------------------------------------------
$ws.Range("F19").Select()
